$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldGuid = "4d9fa03c-8a42-47fc-a0a0-0589e43c8fdc"
$newGuid = "a36eee9f-1c81-4cf0-b9a7-7ed205148f7a"

$oldZhHash = "3030feba0574becae681beb0feabeb69633e6d64"
$newZhHash = "277b1a8ef1fdd8287878bb058a9a454b8ec1dad1"

# --- Overview sheet ---
$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
$overview.Range("G2").Value = "2016-08-25 17:00:10"

# --- zh-cn sheet ---
$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-25 16:59:57"

# --- de-de sheet ---
$dede.Range("A2").Value = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$dede.Range("H2").Value = "2016-08-25 17:00:10"

# --- Update hyperlink display text to match new file name ---
$overview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$zhcn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$dede.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
